# Rename newlink.Stream -> ServerStream; fix comment text that mentions it.
#
# The "Protocol" worksheet (sheet1) has a cell (F5) whose shared-string
# value documents StreamCreatedFrame and currently reads "...that a Stream
# has been successfully created...". The commit renames that concept to
# "ServerStream", so the cell text is updated accordingly. Re-writing the
# cell's value also naturally moves the affected shared-string entry to the
# end of the shared-string table (its old slot is reclaimed) exactly like
# a manual Excel edit would.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F5").Value = "Sent to indicate that a ServerStream has been successfully created. This is the first frame sent over *every* sucessfully-authenticated transport with ``requestNewStream``, so it may be sent over more than one in transport. This allows the client to know that it can now send smaller HelloFrames without ``requestNewStream`` and ``credentialsData``."

# The saved sheet view's selection also changed to a "select all cells"
# state (sqref spans the whole grid), so reproduce that via Ctrl+A-style
# selection of the sheet.
$ws.Cells.Select() | Out-Null
